# Apply the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.986.23'
$ws.Range('E2').Value = '  +2.77%  '
$ws.Range('D3').Value = '2.961.90'
$ws.Range('E3').Value = '  +0.99%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.09'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.64'
$ws.Range('E6').Value = '  +1.01%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').Value = '2.960.26'
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.508'
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.25'
$ws.Range('E10').Value = '  +3.34%  '
$ws.Range('E11').Value = '  +6.37%  '
$ws.Range('E12').Value = '  +0.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000240'
$ws.Range('E13').Value = '  +6.33%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.16'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '3.452.86'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '62.889.52'
$ws.Range('E17').Value = '  +2.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.75'
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').Value = '2.961.21'
$ws.Range('E19').Value = '  +1.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '444.06'
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.53'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.10'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.56'
$ws.Range('E24').Value = '  -0.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.09'
$ws.Range('E25').Value = '  +0.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.91'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('E27').Value = '  -3.77%  '
$ws.Range('E28').Value = '  +0.03%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.35'
$ws.Range('E29').Value = '  +4.91%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.18'
$ws.Range('E31').Value = '  -2.60%  '
$ws.Range('E32').Value = '  +9.76%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.54'
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.109'
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.997'
$ws.Range('E36').Value = '  -1.56%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.15'
$ws.Range('E37').Value = '  +4.67%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.67'
$ws.Range('E38').Value = '  -0.18%  '
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.56'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.53'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('E42').Value = '  -5.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.282'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '39.82'
$ws.Range('E44').Value = '  -6.51%  '
$ws.Range('D45').Value = '2.720.62'
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '134.06'
$ws.Range('E46').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '360.81'
$ws.Range('E48').Value = '  -2.83%  '
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '22.87'
$ws.Range('E51').Value = '  -4.45%  '
